$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $rng = $d.Content
    $found = $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Find failed for: $old"
    }
    $rng.Text = $new
}

Replace-Text "You’re invited to join a focus group discussion about your experience with our programme. This interview is part of a study carried out by researchers from the Universities of Cape Town in South Africa and the University of Oxford in the United Kingdom. " "Jy is uitgenooi om aan te sluit by 'n fokusgroepbespreking oor jou ervaring met ons program. Hierdie onderhoud is deel van 'n studie wat uitgevoer word deur navorsers van die Universiteit van Kaapstad in Suid-Afrika en die Universiteit van Oxford in die Verenigde Koninkryk. "
Replace-Text "Before you decide if you’d like to be interviewed, it’s important for you to know why we’re doing this research and what participating in it would involve. All the information you might need is explained below but if you have any questions about your participation or our study, please email the study team at " "Voordat jy besluit of jy aan die onderhoud wil deelneem, is dit belangrik om te weet waarom ons hierdie navorsing doen en wat deelname behels. Alle inligting wat jy mag nodig hê, word hieronder verduidelik, maar as jy enige vrae het oor jou deelname of die studie, kan jy die navorsingspan per e-pos kontak by "
Replace-Text " or message us on WhatsApp at +27 XX XXX XXXX. We’re here to help you!" " of 'n WhatsApp boodskap stuur na +27 XX XXX XXXX. Ons is hier om jou te help!"
Replace-Text "A focus group is a group discussion for research. As part of this study, you will be placed in a group of 6 – 8 individuals. The other members of the group will all be fellow parents and caregivers who have also been working through the parenting support on the chatbot. A moderator will ask the group questions that will lead to discussion. There might also be a note-taker/interpreter in the room with the group. Both the moderator and the note-taker/interpreter are part of the study team. " "'n Fokusgroep is 'n groepsbespreking wat vir navorsing gebruik word. As deel van hierdie studie sal jy in 'n groep van 6 – 8 individue geplaas word. Die ander lede van die groep sal ook ouers en versorgers wees wat dieselfde ouerskap ondersteuning via die geselsbot gebruik het. 'n Moderator sal vrae aan die groep stel wat tot bespreking sal lei. Daar mag ook 'n notuleerder/tolk in die kamer wees saam met die groep. Beide die moderator en die notuleerder/tolk maak deel uit van die navorsingspan. "
Replace-Text "You’ve been invited to a group discussion (along with other parents/caregivers) with a member of our research team because you’re part of our study. We would love to hear about your experience with the ParentText programme. To be interviewed, you need to agree to take part. " "Jy is uitgenooi om aan 'n groepsbespreking deel te neem (saam met ander ouers/versorgers) met 'n lid van ons navorsingspan omdat jy deel is van ons studie. Ons sal graag jou ervaring met die ParentText-program wil hoor. Om aan die onderhoud deel te neem, moet jy instem om deel te neem. "
Replace-Text "No, it's up to you if you want to join or not. If you don't want to be interviewed, there will be no implications to you or your family. If you do choose to participate in the groups but don't want to answer some of the questions, you can stop at any time by telling your interviewer or just to stop responding in the group." "Nee, dit is jou keuse of jy aansluit of nie. As jy nie wil ondervra word nie, sal daar geen gevolge vir jou of jou familie wees nie. As jy kies om aan die groep deel te neem, maar nie al die vrae wil beantwoord nie, kan jy enige tyd ophou deur die onderhoudvoerder te vertel of net te stop om in die groep te antwoord."
Replace-Text "If you decide you’d like to be interviewed, you’ll need to agree verbally to the consent questions below which the person interviewing you will ask you. The interviews will happen in person and last approximately 1-1.5 hours. The discussion will be conducted in a private space arranged by the research team. " "As jy besluit dat jy graag ondervra wil word, sal jy verbaal moet instem tot die toestemmingsvrae hieronder wat die persoon wat jou ondervra aan jou sal stel. Die onderhoude sal in persoon plaasvind en sal ongeveer 1-1.5 ure neem. Die bespreking sal in 'n privaat ruimte gehou word wat deur die navorsingspan gereël is. "
Replace-Text "During the interview, a member of the research team will ask you some questions about your thoughts and experiences using the chatbot. We want to see if parents like using the chatbot. We also want to know if they're happy with the messages and if using the chatbot changes how they take care of their kids. You will have the right to skip questions you do not want to answer. There are also no right or wrong answers because your whole experience is important to us. " "Tydens die onderhoud sal 'n lid van die navorsingspan jou vrae vra oor jou gedagtes en ervarings met die geselsbot. Ons wil sien of ouers daarvan hou om die geselsbot te gebruik. Ons wil ook weet of hulle tevrede is met die boodskappe en of die gebruik van die geselsbot verander hoe hulle hul kinders versorg. Jy het die reg om enige vrae wat jy nie wil antwoord nie, oor te slaan. Daar is ook geen regte of verkeerde antwoorde nie, want jou hele ervaring is vir ons belangrik. "
Replace-Text "To protect your personal information (including your real name, contact details, and any other information that can identify you), we will give you a participant number, and you can choose a name you want us to call you during the interview. Please also do not refer to any other third parties by name during the interview, without their permission, so that we can protect their personal information too. " "Om jou persoonlike inligting (insluitend jou regte naam, kontakbesonderhede, en enige ander inligting wat jou kan identifiseer) te beskerm, sal ons vir jou 'n deelnemernommer gee, en jy kan 'n naam kies waarmee ons jou tydens die onderhoud kan aanspreek. Moet asseblief nie derde partye se name, sonder hulle toestemming, noem gedurende die onderhoud nie, sodat ons ook hulle persoonlike inligting kan beskerm. "
Replace-Text "We will record the interview to help us remember the discussion and later write down what was said. You will be given a number instead of your name being used so that any information you share in your group interview will not be able to be linked to you by anyone besides the research team. We will delete any personal information we collect from you at the end of the study and, after transcribing your interview, change any data which might lead to identification at the point of transcription. We may use an artificial intelligence (AI) software, Microsoft Transcriber, to transcribe the interviews at first, and then we will check/review these transcriptions. This AI-generated information will be processed and stored securely on password-protected University of Cape Town servers, and in accordance with POPIA. Only authorised members of the research team will be able to access it, and this data will be owned by the Global Parenting Initiative at the University of Cape Town." "Ons sal die onderhoud opneem om die gesprek te onthou en later neer te skryf wat gesê is. Jy sal 'n nommer gegee word in plaas van jou naam, sodat enige inligting wat jy in die groepsonderhoud deel, nie aan jou gekoppel kan word nie, behalwe deur die navorsingspan. Ons sal enige persoonlike inligting wat ons van jou insamel, aan die einde van die studie verwyder, en na die transkripsie van jou onderhoud sal ons enige data wat tot identifisering kan lei, verander op die punt van transkripsie. Ons mag 'n kunsmatige intelligensie (KI) sagteware, Microsoft Transcriber, gebruik om die onderhoude aanvanklik te transkribeer, en daarna sal ons hierdie transkripsies nagaan/hersien. Hierdie KI-geïntegreerde inligting sal verwerk en veilig op wagwoordbeskermde Universiteit van Kaapstad bedieners gestoor word, en in ooreenstemming met POPIA. Slegs gemagtigde lede van die navorsingspan sal toegang hê daartoe, en hierdie data sal deur die Global Parenting Initiative aan die Universiteit van Kaapstad besit word."
Replace-Text "We ask you to respect the other people in the group, and not to discuss what is said by others, outside of the group discussion. We will make sure that our reports are written so that no-one can identify you from the report. Please remember, though, that we can only guarantee this for the research team." "Ons vra dat jy die ander mense in die groep respekteer, en nie wat deur ander gesê word, buite die groep bespreek nie. Ons sal seker maak dat ons verslae op so 'n manier geskryf word dat niemand jou uit die verslag kan identifiseer nie. Onthou asseblief, ons kan dit egter net vir die navorsingspan waarborg."
Replace-Text "We only collect what’s needed for the study and store it securely. Your information, like your consent form and interview recording, and any information you provide via email or WhatsApp, will be kept safe on secure servers at the University of Cape Town. " "Ons versamel slegs wat nodig is vir die studie en stoor dit veilig. Jou inligting, soos jou toestemmingsvorm en onderhoudopname, en enige inligting wat jy via e-pos of WhatsApp verskaf, sal veilig op bedieners by die Universiteit van Kaapstad gestoor word. "
Replace-Text "Interview recordings will be deleted after we have written our notes. Any details that identify you will be kept separate and only authorised staff can access them. All data will be kept for five years after the study, but personal information will be deleted when the study ends. " "Onderhoudopnames sal verwyder word nadat ons, ons notas geskryf het. Enige besonderhede wat jou kan identifiseer, sal afsonderlik gehou word en slegs gemagtigde personeel sal toegang hê. Alle data sal vyf jaar na die studie bewaar word, maar persoonlike inligting sal verwyder word wanneer die studie eindig. "
Replace-Text "Ethics committees and monitors may check the information. Your information will stay private unless the law says otherwise. After the study, we may share the information with other researchers but without your details. You have the right to see, correct, or ask us to delete your personal information." "Etiekkomitees en moniteerders mag die inligting nagaan. Jou inligting sal privaat bly tensy die wet anders bepaal. Na die studie mag ons die inligting met ander navorsers deel, maar sonder jou besonderhede. Jy het die reg om jou persoonlike inligting te sien, reg te stel of ons te vra om jou persoonlike inligting te verwyder."
Replace-Text "You have the right to request access to your data, to correct any mistakes in your data, and to request us to delete it or transfer it somewhere else. Please email the study team before [*date to be determined] if you would like to do any of these." "Jy het die reg om toegang tot jou data aan te vra, enige foute in jou data reg te stel en te vra dat ons dit verwyder of dit na 'n ander plek oordra. Stuur asseblief 'n e-pos aan die studiespan voor [*datum moet nog bepaal word] indien jy enige van hierdie aksies wil uitvoer."
Replace-Text "Your participation and what you tell us will help us understand how to support families like yours. We plan to share the results in reports and at conferences so others can learn from this study too." "Jou deelname en wat jy ons vertel, sal ons help om te verstaan hoe ons gesinne soos joune beter kan ondersteun. Ons beplan om die resultate in verslae en op kongresse te deel sodat ander ook van hierdie studie kan leer."
Replace-Text "Are there any risks in being interviewed?   " "Is daar enige risiko's verbonde aan die onderhoud?   "
Replace-Text "We don’t expect any risks to you if you are interviewed. If any questions make you uncomfortable, you don’t have to answer them. If you become upset during the interview, you can let your interviewer know. Remember, you can stop participating anytime without giving a reason. We care about your well-being." "Ons verwag nie enige risiko's vir jou as jy ondervra word nie. As enige vrae jou ongemaklik laat voel, hoef jy dit nie te antwoord nie. As jy tydens die onderhoud ontsteld raak, kan jy jou onderhoudvoerder laat weet. Onthou, jy kan enige tyd ophou deelneem sonder om 'n rede te gee. Ons gee om vir jou welstand."
Replace-Text "We also want to make sure you're safe. If we notice that you or your family are in serious danger, we might refer you for support or could need to ask for help from other places outside of this study, like social or medical services." "Ons wil ook seker maak jy is veilig. As ons opmerk dat jy of jou gesin in ernstige gevaar is, mag ons jou vir ondersteuning verwys of kan vir hulp vra van ander plekke buite hierdie studie, soos maatskaplike of mediese dienste."
Replace-Text "The University Cape Town makes sure your personal information is used safely and correctly, just for research. The study follows data protection laws like GDPR (General Data Protection Regulation) in the UK and POPIA (Protection of Personal Information Act) in South Africa. Any data that is transferred across borders will comply with POPIA. " "Die Universiteit van Kaapstad sorg dat jou persoonlike inligting veilig en korrek gebruik word, net vir navorsingsdoeleindes. Die studie volg databeskermingswette soos die GDPR (General Data Protection Regulation) in die VK en POPIA (Wet op die Beskerming van Persoonlike Inligting) in Suid-Afrika. Enige data wat oor grense heen oorgedra word, sal aan POPIA voldoen. "
Replace-Text "[Once the ethics has been approved this will read as follows: This study has received approval from the University of Cape Town’s Centre for Social Science Research Ethics Committee and University of Cape Town’s Faculty of Health Sciences Human Research Ethics Committee. The study has also been approved by the Western Cape Department of Health and Wellness a" "[Wanneer die etiek goedgekeur is, sal dit as volg lees: Hierdie studie het goedkeuring van die Universiteit van Kaapstad se Sentrum vir Sosiale Wetenskap Navorsingsetiekkomitee en die Universiteit van Kaapstad se Fakulteit van Gesondheidswetenskappe se Etiekkomitee vir Menslike Navorsing ontvang. Die studie is ook goedgekeur deur die Wes-Kaapse Departement van Gesondheid en Welstand e"
Replace-Text "nd Department of Social Development," "n die Departement van Maatskaplike Ontwikkeling,"
Replace-Text " and City of Cape Town’s City health.]" " en die Stad Kaapstad se Stadsgesondheid.]"
Replace-Text "I had time to think about the information and ask questions. I am happy with the answers which I got. " "Ek het tyd gehad om oor die inligting te dink en vrae te vra. Ek is tevrede met die antwoorde wat ek ontvang het. "
Replace-Text "I know I can say yes or no to being in the study. If I say yes, I can stop any time before the [*date to still be determined] without saying why, and nothing bad will happen." "Ek weet ek kan ja of nee sê om aan die studie deel te neem. As ek ja sê, kan ek enige tyd voor die [*datum wat nog bepaal moet word] stop sonder om te sê hoekom, en niks sleg sal gebeur nie."
Replace-Text "I am okay with the interview being recorded. I know the recordings will be used for research." "Ek is reg daarmee dat die onderhoud opgeneem word. Ek weet die opnames sal gebruik word vir navorsing."
Replace-Text "I will respect others in the group and will not discuss what is said by others outside of the group. I understand that the other participants in the group have to keep this same promise." "Ek sal ander in die groep respekteer en nie buite die groep bespreek wat ander gesê het nie. Ek verstaan dat die ander deelnemers in die groep dieselfde belofte moet maak."

Write-Output "done"